# Facilitators guidelines - Moebius.docx : English -> Swahili (Kenya) translation pass.
#
# Each English label/phrase in the table is replaced with its Swahili
# translation. "General VMC Video Introduction" is replaced before the
# shorter "Video Introduction" because the former contains the latter as a
# substring; doing the longer replacement first avoids corrupting it.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "Video Title" "Kichwa cha Video"
Replace-Text "Topic" "Mada"
Replace-Text "Geometry" "Jiometri"
Replace-Text "Aim(s)" "Malengo"
Replace-Text "Length" "Urefu"
Replace-Text "Camp Location" "Mahali pa Kambi"
Replace-Text "Facilitators" "Wawezeshaji"
Replace-Text "N. of students" "N. ya wanafunzi"
Replace-Text "Date" "Tarehe"
Replace-Text "Resources" "Rasilimali"
Replace-Text "needed" "inahitajika"
Replace-Text "Preparations" "Maandalizi"
Replace-Text "Video time" "Muda wa video"
Replace-Text "What facilitator does" "Mwezeshaji anafanya nini"
Replace-Text "What learners do" "Wanachofanya wanafunzi"
Replace-Text "General VMC Video Introduction" "Utangulizi Mkuu wa Video ya VMC"
Replace-Text "Video Introduction" "Utangulizi wa Video"
Replace-Text "Introduction of the first experiment" "Utangulizi wa jaribio la kwanza"
Replace-Text "Assist the process, provoke thoughts" "Kusaidia mchakato, kuchochea mawazo"

# Default document language: Swahili (Tanzania) -> Swahili (Kenya).
# wdStyleNormal (-1) is the base style every paragraph/run in this document
# inherits its language from, so retargeting it is the COM-level equivalent
# of flipping the document's default w:lang.
$normalStyle = $d.Styles.Item(-1)
$normalStyle.LanguageID = "sw-KE"
